$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific cells that would otherwise be auto-converted to numbers
# by Excel's type inference to remain plain text, matching the source data
# (the workbook stores all Price/Volume cells as text).
$forceTextCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D33", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D46", "D48")
foreach ($c in $forceTextCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "64.479.39"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "3.152.44"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "594.33"
$ws.Range("D6").Value = "146.95"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.150.83"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "5.91"
$ws.Range("E11").Value = "  +3.21%  "
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").Value = "37.45"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").Value = "3.678.26"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "7.31"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D18").Value = "64.244.01"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "3.156.11"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").Value = "470.63"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").Value = "0.738"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "7.58"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "13.23"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "81.71"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "2.33"
$ws.Range("E26").Value = "  +5.88%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +7.47%  "
$ws.Range("D29").Value = "7.48"
$ws.Range("E29").Value = "  +8.64%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "27.61"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = ("0.0{0}0849" -f [char]0x2083)
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "6.26"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").Value = "51.94"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").Value = "458.76"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("D42").Value = "9.30"
$ws.Range("E42").Value = "  +5.89%  "
$ws.Range("D43").Value = "0.297"
$ws.Range("E43").Value = "  +6.60%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").Value = "2.941.96"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "40.55"
$ws.Range("E46").Value = "  +12.44%  "
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D48").Value = "129.35"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("E51").Value = "  -0.21%  "
